# Add files via upload
# Extends the MAR-2021 sheet (sheet7) with rows 21-32 (20-31 Mar 2021),
# reusing the same cell-style patterns already present in rows 2-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAR-2021")

# ---- Row 21 : 20 Mar 2021 - B2B api testing / B2B app (Completed) ----
$ws.Range("A20").Copy($ws.Range("A21"))
$ws.Range("A21").Value = 20
$ws.Range("B20").Copy($ws.Range("B21"))
$ws.Range("B21").Value = 44275
$ws.Range("C7").Copy($ws.Range("C21"))
$ws.Range("D20").Copy($ws.Range("D21"))
$ws.Range("D21").Value = "B2B api testing"
$ws.Range("C21").Value = "B2B app"
$ws.Range("E20").Copy($ws.Range("E21"))
$ws.Range("E21").Value = 1
$ws.Range("F20").Copy($ws.Range("F21"))
$ws.Range("G20").Copy($ws.Range("G21"))

# ---- Row 22 : 21 Mar 2021 - Week off ----
$ws.Range("A20").Copy($ws.Range("A22"))
$ws.Range("A22").Value = 21
$ws.Range("B20").Copy($ws.Range("B22"))
$ws.Range("B22").Value = 44276
$ws.Range("C7").Copy($ws.Range("C22"))
$ws.Range("D7").Copy($ws.Range("D22"))
$ws.Range("E7").Copy($ws.Range("E22"))
$ws.Range("F7").Copy($ws.Range("F22"))
$ws.Range("G7").Copy($ws.Range("G22"))

# ---- Row 23 : 22 Mar 2021 - nMVAR and Muji store (Completed) ----
$ws.Range("A20").Copy($ws.Range("A23"))
$ws.Range("A23").Value = 22
$ws.Range("B20").Copy($ws.Range("B23"))
$ws.Range("B23").Value = 44277
$ws.Range("C20").Copy($ws.Range("C23"))
$ws.Range("C23").Value = "nMVAR and Muji store"
$ws.Range("D21").Copy($ws.Range("D23"))
$ws.Range("D23").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. `nRegression testing and Retesting on nMVAR_Report, nMVAR_TSS, nMVAR_Invc and nMVAR_QA application`nRegression testing and Retesting on Muji store application`nCross browser testing on Muji store application"
$ws.Range("E20").Copy($ws.Range("E23"))
$ws.Range("E23").Value = 1
$ws.Range("F20").Copy($ws.Range("F23"))
$ws.Range("G20").Copy($ws.Range("G23"))
$ws.Rows.Item(23).RowHeight = 75

# ---- Rows 24-32 : 23 Mar 2021 - 31 Mar 2021 (blank placeholder rows) ----
$blankDates = @(44278, 44279, 44280, 44281, 44282, 44283, 44284, 44285, 44286)
$blankNos = @(23, 24, 25, 26, 27, 28, 29, 30, 31)
for ($i = 0; $i -lt $blankDates.Length; $i++) {
    $r = 24 + $i
    $ws.Range("A20").Copy($ws.Range("A$r"))
    $ws.Range("A$r").Value = $blankNos[$i]
    $ws.Range("B20").Copy($ws.Range("B$r"))
    $ws.Range("B$r").Value = $blankDates[$i]
    $ws.Range("C7").Copy($ws.Range("C$r"))
    $ws.Range("D21").Copy($ws.Range("D$r"))
    $ws.Range("D$r").ClearContents()
    $ws.Range("E7").Copy($ws.Range("E$r"))
    $ws.Range("F7").Copy($ws.Range("F$r"))
    $ws.Range("G20").Copy($ws.Range("G$r"))
}

# ---- Update the visible selection to match the saved view (D29) ----
$ws.Range("D29").Select()
